$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 129
$ws.Cells.Item(129, 8).Value = 989.8333
$ws.Cells.Item(129, 9).Value = 0
$ws.Cells.Item(129, 10).Value = 989.8333
$ws.Cells.Item(129, 11).Value = 0
$ws.Cells.Item(129, 12).Value = 2969.4999
$ws.Cells.Item(129, 13).ClearContents()
$ws.Cells.Item(129, 14).Value = -12969.4999

$ws = $wb.Worksheets.Item("ARM")
# ARM row 45
$ws.Cells.Item(45, 8).Value = 1754.7354
$ws.Cells.Item(45, 9).Value = 1645.2069
$ws.Cells.Item(45, 10).Value = 2390
$ws.Cells.Item(45, 11).Value = 1645.2069
$ws.Cells.Item(45, 12).Value = 2390
$ws.Cells.Item(45, 13).Value = -1268.2069
$ws.Cells.Item(45, 14).Value = -3144

# ARM row 61
$ws.Cells.Item(61, 8).Value = 6501.8604
$ws.Cells.Item(61, 9).Value = 3886.516
$ws.Cells.Item(61, 10).Value = 13258.167
$ws.Cells.Item(61, 11).Value = 3886.516
$ws.Cells.Item(61, 12).Value = 13258.167
$ws.Cells.Item(61, 13).Value = -3674.516
$ws.Cells.Item(61, 14).Value = -13682.167

# ARM row 113
$ws.Cells.Item(113, 8).Value = 79800
$ws.Cells.Item(113, 10).Value = 79800
$ws.Cells.Item(113, 12).Value = 79800
$ws.Cells.Item(113, 14).Value = -88478

# ARM row 122
$ws.Cells.Item(122, 8).Value = 5001760
$ws.Cells.Item(122, 9).Value = 1839.1
$ws.Cells.Item(122, 10).Value = 25001444
$ws.Cells.Item(122, 11).Value = 5517.299999999999
$ws.Cells.Item(122, 12).Value = 75004332
$ws.Cells.Item(122, 13).Value = -3067.299999999999
$ws.Cells.Item(122, 14).Value = -75009232

# ARM row 128
$ws.Cells.Item(128, 8).Value = 56150
$ws.Cells.Item(128, 10).Value = 56150
$ws.Cells.Item(128, 12).Value = 56150
$ws.Cells.Item(128, 14).Value = -66110

# ARM row 132
$ws.Cells.Item(132, 8).Value = 1998.283
$ws.Cells.Item(132, 9).Value = 1565.9773
$ws.Cells.Item(132, 10).Value = 4111.778
$ws.Cells.Item(132, 11).Value = 4697.9319
$ws.Cells.Item(132, 12).Value = 12335.334
$ws.Cells.Item(132, 13).Value = -2167.9319
$ws.Cells.Item(132, 14).Value = -17395.334

# ARM row 136
$ws.Cells.Item(136, 8).Value = 6501.8604
$ws.Cells.Item(136, 9).Value = 3886.516
$ws.Cells.Item(136, 10).Value = 13258.167
$ws.Cells.Item(136, 11).Value = 11659.548
$ws.Cells.Item(136, 12).Value = 39774.501
$ws.Cells.Item(136, 13).Value = -9109.548000000001
$ws.Cells.Item(136, 14).Value = -44874.501

$ws = $wb.Worksheets.Item("BSM")
# BSM row 55
$ws.Cells.Item(55, 8).Value = 0
$ws.Cells.Item(55, 10).Value = 0
$ws.Cells.Item(55, 12).Value = 0
$ws.Cells.Item(55, 14).ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# CRP row 122
$ws.Cells.Item(122, 8).Value = 16561.223
$ws.Cells.Item(122, 10).Value = 30257
$ws.Cells.Item(122, 12).Value = 90771
$ws.Cells.Item(122, 14).Value = -95671

$ws = $wb.Worksheets.Item("CUL")
# CUL row 63
$ws.Cells.Item(63, 8).Value = 3085.1667
$ws.Cells.Item(63, 9).Value = 1304
$ws.Cells.Item(63, 10).Value = 4866.3335
$ws.Cells.Item(63, 11).Value = 3912
$ws.Cells.Item(63, 12).Value = 14599.0005
$ws.Cells.Item(63, 13).Value = -3163
$ws.Cells.Item(63, 14).Value = -16097.0005

# CUL row 66
$ws.Cells.Item(66, 8).Value = 3085.1667
$ws.Cells.Item(66, 9).Value = 1304
$ws.Cells.Item(66, 10).Value = 4866.3335
$ws.Cells.Item(66, 11).Value = 11736
$ws.Cells.Item(66, 12).Value = 43797.0015
$ws.Cells.Item(66, 13).Value = -7992
$ws.Cells.Item(66, 14).Value = -51285.0015

# CUL row 68
$ws.Cells.Item(68, 8).Value = 2859.6626
$ws.Cells.Item(68, 9).Value = 1508.9756
$ws.Cells.Item(68, 10).Value = 4279.615
$ws.Cells.Item(68, 11).Value = 4526.9268
$ws.Cells.Item(68, 12).Value = 12838.845
$ws.Cells.Item(68, 13).Value = -3715.9268
$ws.Cells.Item(68, 14).Value = -14460.845

# CUL row 71
$ws.Cells.Item(71, 8).Value = 2859.6626
$ws.Cells.Item(71, 9).Value = 1508.9756
$ws.Cells.Item(71, 10).Value = 4279.615
$ws.Cells.Item(71, 11).Value = 13580.7804
$ws.Cells.Item(71, 12).Value = 38516.535
$ws.Cells.Item(71, 13).Value = -9524.7804
$ws.Cells.Item(71, 14).Value = -46628.535

# CUL row 81
$ws.Cells.Item(81, 8).Value = 4081.875
$ws.Cells.Item(81, 9).Value = 1525
$ws.Cells.Item(81, 10).Value = 4934.1665
$ws.Cells.Item(81, 11).Value = 4575
$ws.Cells.Item(81, 12).Value = 14802.4995
$ws.Cells.Item(81, 13).Value = -3452
$ws.Cells.Item(81, 14).Value = -17048.4995

# CUL row 84
$ws.Cells.Item(84, 8).Value = 4081.875
$ws.Cells.Item(84, 9).Value = 1525
$ws.Cells.Item(84, 10).Value = 4934.1665
$ws.Cells.Item(84, 11).Value = 13725
$ws.Cells.Item(84, 12).Value = 44407.4985
$ws.Cells.Item(84, 13).Value = -8109
$ws.Cells.Item(84, 14).Value = -55639.4985

# CUL row 131
$ws.Cells.Item(131, 8).Value = 1072.8909
$ws.Cells.Item(131, 9).Value = 1173.2778
$ws.Cells.Item(131, 10).Value = 1024.0541
$ws.Cells.Item(131, 11).Value = 3519.8334
$ws.Cells.Item(131, 12).Value = 3072.1623
$ws.Cells.Item(131, 13).Value = 1520.1666
$ws.Cells.Item(131, 14).Value = -13152.1623

# CUL row 137
$ws.Cells.Item(137, 8).Value = 43900.668
$ws.Cells.Item(137, 9).Value = 1476
$ws.Cells.Item(137, 10).Value = 128750
$ws.Cells.Item(137, 11).Value = 4428
$ws.Cells.Item(137, 12).Value = 386250
$ws.Cells.Item(137, 13).Value = 672
$ws.Cells.Item(137, 14).Value = -396450

$ws = $wb.Worksheets.Item("GSM")
# GSM row 132
$ws.Cells.Item(132, 8).Value = 7404.4443
$ws.Cells.Item(132, 9).Value = 1575.091
$ws.Cells.Item(132, 10).Value = 16564.857
$ws.Cells.Item(132, 11).Value = 4725.272999999999
$ws.Cells.Item(132, 12).Value = 49694.571
$ws.Cells.Item(132, 13).Value = -2195.272999999999
$ws.Cells.Item(132, 14).Value = -54754.571

$ws = $wb.Worksheets.Item("LTW")
# LTW row 7
$ws.Cells.Item(7, 8).Value = 2311.1936
$ws.Cells.Item(7, 9).Value = 2517.0557
$ws.Cells.Item(7, 10).Value = 2026.1538
$ws.Cells.Item(7, 11).Value = 2517.0557
$ws.Cells.Item(7, 12).Value = 2026.1538
$ws.Cells.Item(7, 13).Value = -2405.0557
$ws.Cells.Item(7, 14).Value = -2250.1538

# LTW row 16
$ws.Cells.Item(16, 8).Value = 1195
$ws.Cells.Item(16, 9).Value = 933.3333
$ws.Cells.Item(16, 10).Value = 1980
$ws.Cells.Item(16, 11).Value = 933.3333
$ws.Cells.Item(16, 12).Value = 1980
$ws.Cells.Item(16, 13).Value = -763.3333
$ws.Cells.Item(16, 14).Value = -2320

# LTW row 40
$ws.Cells.Item(40, 8).Value = 3416.1304
$ws.Cells.Item(40, 9).Value = 3285.7144
$ws.Cells.Item(40, 10).Value = 3619
$ws.Cells.Item(40, 11).Value = 3285.7144
$ws.Cells.Item(40, 12).Value = 3619
$ws.Cells.Item(40, 13).Value = -3149.7144
$ws.Cells.Item(40, 14).Value = -3891

# LTW row 122
$ws.Cells.Item(122, 8).Value = 7321.5713
$ws.Cells.Item(122, 9).Value = 6869.522
$ws.Cells.Item(122, 10).Value = 9401
$ws.Cells.Item(122, 11).Value = 20608.566
$ws.Cells.Item(122, 12).Value = 28203
$ws.Cells.Item(122, 13).Value = -18158.566
$ws.Cells.Item(122, 14).Value = -33103

# LTW row 126
$ws.Cells.Item(126, 8).Value = 2311.1936
$ws.Cells.Item(126, 9).Value = 2517.0557
$ws.Cells.Item(126, 10).Value = 2026.1538
$ws.Cells.Item(126, 11).Value = 7551.1671
$ws.Cells.Item(126, 12).Value = 6078.4614
$ws.Cells.Item(126, 13).Value = -5081.1671
$ws.Cells.Item(126, 14).Value = -11018.4614

$ws = $wb.Worksheets.Item("WVR")
# WVR row 122
$ws.Cells.Item(122, 8).Value = 2151.125
$ws.Cells.Item(122, 9).Value = 1900.6666
$ws.Cells.Item(122, 10).Value = 2902.5
$ws.Cells.Item(122, 11).Value = 5701.9998
$ws.Cells.Item(122, 12).Value = 8707.5
$ws.Cells.Item(122, 13).Value = -3251.9998
$ws.Cells.Item(122, 14).Value = -13607.5

# WVR row 132
$ws.Cells.Item(132, 8).Value = 1744.9546
$ws.Cells.Item(132, 9).Value = 1744.5161
$ws.Cells.Item(132, 10).Value = 1746
$ws.Cells.Item(132, 11).Value = 5233.5483
$ws.Cells.Item(132, 12).Value = 5238
$ws.Cells.Item(132, 13).Value = -2703.5483
$ws.Cells.Item(132, 14).Value = -10298
